# Update the "Förändrad" (changed) date column C for rows 2-8:
# the underlying date serial value moves from 45183 (2023-09-14)
# to 45184 (2023-09-15), while keeping existing number formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value = 45184
    }
}
